$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New state-transition UUIDs for "ingest_session_issue_id" (column A), one per
# data row (rows 2-49), issued as part of introducing state transitions logging.
$newIssueIds = @(
    "6cb3bee1-de9a-4154-99a6-ef593d8ed12d",
    "d44d87de-4b8d-4a70-841e-560e7c9f1e68",
    "93fbfaea-1a23-4aa1-b225-e34b1e1d5c9c",
    "73aee1e3-15de-49c8-b241-152dfb69cd6f",
    "19f9afdb-bb49-4c2b-a96a-61696eb0b2bd",
    "7a8328d7-0f44-405d-a2ff-9bfcefc606d7",
    "a5b093e0-4b82-47a7-b638-6a5dd53158c0",
    "4ad6f070-d21d-4a76-9c72-01d5cffd4e35",
    "6bf2021a-770e-496b-b513-a1b7bd3fd29d",
    "51a1dbda-64b3-4697-9a8c-077973981313",
    "12258278-c453-4dac-9e95-d121b0c83cdf",
    "8f74cf72-368d-41b1-a069-d8076087fcf0",
    "ea231125-6854-446e-aefa-4c6d8572ae80",
    "f7b488fa-38b5-4ebe-a85b-747c38a4c185",
    "40077d9f-b711-4699-8f42-6312e3b45597",
    "17969a5c-9cc1-45c8-9445-c1112c87809a",
    "88367295-7e06-419b-958c-6f639eea4096",
    "32045c01-7ce0-490c-a94f-5123e46793ac",
    "1678fa3e-5c7a-4854-9edc-9e5bfafed84f",
    "5debc646-d08d-4286-bc74-53d398f18891",
    "1a276d48-d7cf-47e4-9bf3-f8bf1340ea99",
    "1ada402b-a020-4b1c-b9dc-1197b283c925",
    "f21d9cee-ef24-412b-8aeb-f446b3c1ee67",
    "aebeb812-9558-4b49-b37f-e4795f292354",
    "5266de63-fc55-4289-9d9b-f927652071d9",
    "c0edb199-c40a-4139-be8f-249297399436",
    "f8d5b20e-f651-490c-af48-fa3cd1fbee04",
    "1e78492f-b48b-43aa-9d61-0e986dd547bd",
    "36778307-05a1-42e5-83ec-b60de93ea445",
    "05e1a9b8-f95b-4d36-9fa2-4a8e87ee1b5d",
    "eaa3cd89-d225-4d5d-8b4f-ec05f0e6266f",
    "dfc46fb6-f998-45f3-9449-972116bcc150",
    "1d297a2a-ded9-4f17-9e29-d60a8cc9a36b",
    "739077ce-49d8-4b79-be39-3c4ecbf9d795",
    "f99e260d-3f37-4f6b-a48a-ffa6c48e6c0f",
    "4fbb9a18-d484-4eab-adbb-fa66f4e2433c",
    "6ba59847-caca-479f-8860-9f78f41a45ea",
    "3f0eacb6-bd9c-4adf-8335-22581491630d",
    "9dadc21c-134d-426a-868a-7d89d9c136af",
    "d3872d39-a048-4cf2-8474-cb58505268d9",
    "dcf58158-5a39-4f82-9e69-a689c707d794",
    "8b089552-529e-429f-9759-d35a8d997641",
    "2ab96894-50c8-4f14-9b04-66fc17aef1a5",
    "b40e728f-a8ed-47aa-809b-520540b53112",
    "3d82f5d1-3778-4d1d-87ea-6b2f2ad3ddc2",
    "3288676a-6775-46c5-9cac-4c7067e7d7b2",
    "c747bd11-355a-4612-87e0-7cff07866ad8",
    "b16dee1c-c9e3-4d1d-b398-7325df3e95a4"
)

for ($i = 0; $i -lt $newIssueIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newIssueIds[$i]
}

# New "session_entry_id" (column C) values for the two logging sessions captured
# in this diagnostics sheet.
for ($row = 2; $row -le 22; $row++) {
    $ws.Cells.Item($row, 3).Value = "1931dfcc-e8fc-597d-b1bc-65b4287e6fdf"
}
for ($row = 23; $row -le 49; $row++) {
    $ws.Cells.Item($row, 3).Value = "8b7c669c-1795-5f6b-8f3a-3e502b74c628"
}

